$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Publication date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws1.Range("B9").Value = "Alvearie Team"

# Old duplicate "Contact" row repurposed into a "Jurisdiction" row
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The second (duplicate) "Contact" / "No display for ContactDetail" row is removed entirely,
# shifting everything below it up by one row (and dropping the trailing "Context" row).
$ws1.Rows.Item(11).Delete()

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (the root "Extension" element) gets a real Short/Definition instead of the generic placeholders
$ws2.Range("K2").Value = "Cobra Indicator"
$ws2.Range("L2").Value = "Indicator of Consolidated Omnibus Budget Reconciliation Act(COBRA) continuation for the person"
